$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Append "L" to the kit catalog number in column G (rows 2-17):
#    "NEBNextPoly(A)E7490" -> "NEBNextPoly(A)E7490L"
$ws.Range("G2:G17").Value = "NEBNextPoly(A)E7490L"

# 2. Widen column G to fit the longer text (xlsx <col> width ~27.23)
$ws.Columns.Item(7).ColumnWidth = 26.39666666666667

# 3. Turn the literal FALSE booleans in column I (rows 2-17) into live
#    formulas "=FALSE()" (kept as independent, non-shared formulas so each
#    row carries its own <f> element, matching the per-row diff).
for ($r = 2; $r -le 17; $r++) {
    $ws.Range("I" + $r).Formula = "=FALSE()"
}

# 4. Move the active selection from I2:I17 to G2:G17
$ws.Range("G2:G17").Select()
